$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "In progress" label on row 2 (trailing space added) ---
$ws.Range("A2").Value = "In progress "

# --- 2. Make room for 6 new backlog rows: insert before the old row 8,
#        pushing the previously-blank template rows (10-15) down to (16-21).
$ws.Rows("8:13").Insert()

# --- 3. Clear the stray formatting that used to live on D7 (and that the
#        row-insert copied into D8:D12) so those cells fall back to the
#        workbook default style.
$ws.Range("D7:D12").ClearFormats()

# --- 4. Populate the brand-new task descriptions first, in the same order
#        the original author typed them, so the rebuilt shared-strings
#        table lines up with the canonical file.
$ws.Range("D5").Value = "Implement filter for filtering member orders when updating rental status "
$ws.Range("D7").Value = "Implement SQL queries for gathering all orders made by members for them to view "
$ws.Range("D4").Value = "Implement SQL queries for gathering orders of members"
$ws.Range("D8").Value = "Design page for customer to view their history"
$ws.Range("D11").Value = "Implement SQL queries for gathering information about a single item"
$ws.Range("D12").Value = "Design page for manager to view history"

# --- 5. Fill in the remaining cells for the updated / new task rows ------
# Row 4
$ws.Range("A4").Value = "Complete"
$ws.Range("B4").Value = "Luke"
$ws.Range("C4").Value = "Warehouse employee"
$ws.Range("E4").Value = 43866
$ws.Range("F4").Value = 43873

# Row 5
$ws.Range("A5").Value = "Complete"
$ws.Range("B5").Value = "Luke"
$ws.Range("C5").Value = "Warehouse employee"
$ws.Range("E5").Value = 43866
$ws.Range("F5").Value = 43873

# Row 6
$ws.Range("A6").Value = "Complete"
$ws.Range("B6").Value = "Luke"
$ws.Range("C6").Value = "Warehouse employee"
$ws.Range("D6").Value = "Update rental status of member order"
$ws.Range("E6").Value = 43866
$ws.Range("F6").Value = 43873

# Row 7
$ws.Range("A7").Value = "Complete"
$ws.Range("B7").Value = "Tristen"
$ws.Range("C7").Value = "Member"
$ws.Range("E7").Value = 43869
$ws.Range("F7").Value = 43873

# Row 8 (new)
$ws.Range("A8").Value = "Complete"
$ws.Range("B8").Value = "Tristen"
$ws.Range("C8").Value = "Member"
$ws.Range("E8").Value = 43869
$ws.Range("F8").Value = 43873

# Row 9 (new)
$ws.Range("A9").Value = "Complete"
$ws.Range("B9").Value = "Tristen"
$ws.Range("C9").Value = "Member"
$ws.Range("D9").Value = "View rental history"
$ws.Range("E9").Value = 43866
$ws.Range("F9").Value = 43873

# Row 10 (new)
$ws.Range("A10").Value = "Complete"
$ws.Range("B10").Value = "Luke"
$ws.Range("C10").Value = "Warehouse employee"
$ws.Range("D10").Value = "View detailed list of orders that are outgoing/incoming"
$ws.Range("E10").Value = 43866
$ws.Range("F10").Value = 43873

# Row 11 (new)
$ws.Range("A11").Value = "In progress "
$ws.Range("B11").Value = "Carson"
$ws.Range("C11").Value = "Manager"
$ws.Range("E11").Value = 43871
$ws.Range("F11").Value = 43873

# Row 12 (new)
$ws.Range("A12").Value = "In progress "
$ws.Range("B12").Value = "Carson"
$ws.Range("C12").Value = "Manager"
$ws.Range("E12").Value = 43871
$ws.Range("F12").Value = 43873

# Row 13 (new) - note the distinct font applied to D13 below
$ws.Range("A13").Value = "In Progress "
$ws.Range("B13").Value = "Carson"
$ws.Range("C13").Value = "Manager"
$ws.Range("D13").Value = "view history a single item"
$ws.Range("E13").Value = 43869
$ws.Range("F13").Value = 43873

# D13 carries its own (non-themed) Calibri 11 font in the saved workbook.
$ws.Range("D13").Font.Name = "Calibri"
$ws.Range("D13").Font.Size = 11

# --- 5. Window / view cosmetics ------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 145
